$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cells = @(
    @{addr="C3"; val="0"},
    @{addr="D3"; val="2"},
    @{addr="C4"; val="1"},
    @{addr="D4"; val="1"},
    @{addr="E4"; val="0"},
    @{addr="C5"; val="4"},
    @{addr="D5"; val="7"},
    @{addr="C6"; val="3"},
    @{addr="D6"; val="3"},
    @{addr="C7"; val="8"},
    @{addr="D7"; val="5"},
    @{addr="E7"; val="1"}
)

foreach ($cell in $cells) {
    $rng = $ws.Range($cell.addr)
    $rng.NumberFormat = "@"
    $rng.Value = $cell.val
}
